# Generate Report for Handoff
#
# zh-cn was handed off again (new handoff timestamp), so its status flips
# from "Handed back: in sync with en-US" to "Ready for handoff" across the
# Overview sheet and both locale sheets, and the relevant datetime stamps
# are refreshed. The Status / datetime columns also get narrower now that
# the long "Handed back: in sync with en-US" text is gone.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Datetime stamp refreshes
# Overview "Latest HO Xliff Generate Date" / de-de "Latest Handback DateTime"
$wsOverview.Range("G2").Value = "2016-09-04 13:03:57"
$wsDeDe.Range("H2").Value = "2016-09-04 13:03:57"
# zh-cn "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-04 13:03:53"

# --- Column width shrink now the Status text is shorter
# (ColumnWidth is stored in whole-pixel increments under the hood, so the
# input below is the value that lands on the pixel closest to the target
# serialized width of 17.2159881591797 characters.)
$wsOverview.Range("E1:F1").ColumnWidth = 16.3333333333333
$wsZhCn.Range("C1").ColumnWidth = 16.3333333333333
$wsDeDe.Range("C1").ColumnWidth = 16.3333333333333
